# Update NATMI LR-pairs sheet with new TPM-derived values.
# Rows 2-11 are updated in place, and rows 12-16 are newly added
# (third "sending cluster" = Inflammatory-Mac combined with the 5 target clusters).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lta"
$ws.Range("C2").Value = "Tnfrsf1b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6836903333333334
$ws.Range("H2").Value = 2.051071
$ws.Range("I2").Value = 0.7427603292077943
$ws.Range("J2").Value = 0.7427603292077943
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 18.95273633333333
$ws.Range("N2").Value = 56.858209
$ws.Range("O2").Value = 0.08721078561875104
$ws.Range("P2").Value = 0.08721078561875105
$ws.Range("Q2").Value = 12.95780262131544
$ws.Range("R2").Value = 116.620223591839
$ws.Range("S2").Value = 0.0647767118366539
$ws.Range("T2").Value = 0.0647767118366539

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lta"
$ws.Range("C3").Value = "Tnfrsf1b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6836903333333334
$ws.Range("H3").Value = 2.051071
$ws.Range("I3").Value = 0.7427603292077943
$ws.Range("J3").Value = 0.7427603292077943
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.54486
$ws.Range("N3").Value = 31.63458
$ws.Range("O3").Value = 0.04852204497892696
$ws.Range("P3").Value = 0.04852204497892696
$ws.Range("Q3").Value = 7.209418848353334
$ws.Range("R3").Value = 64.88476963518001
$ws.Range("S3").Value = 0.0360402501023832
$ws.Range("T3").Value = 0.0360402501023832

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lta"
$ws.Range("C4").Value = "Tnfrsf1b"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6836903333333334
$ws.Range("H4").Value = 2.051071
$ws.Range("I4").Value = 0.7427603292077943
$ws.Range("J4").Value = 0.7427603292077943
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 101.4555613333333
$ws.Range("N4").Value = 304.366684
$ws.Range("O4").Value = 0.4668465309523581
$ws.Range("P4").Value = 0.4668465309523581
$ws.Range("Q4").Value = 69.36418654650711
$ws.Range("R4").Value = 624.2776789185641
$ws.Range("S4").Value = 0.3467550830196903
$ws.Range("T4").Value = 0.3467550830196903

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Lta"
$ws.Range("C5").Value = "Tnfrsf1b"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6836903333333334
$ws.Range("H5").Value = 2.051071
$ws.Range("I5").Value = 0.7427603292077943
$ws.Range("J5").Value = 0.7427603292077943
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.410466333333333
$ws.Range("N5").Value = 7.231399000000001
$ws.Range("O5").Value = 0.01109173150200089
$ws.Range("P5").Value = 0.01109173150200089
$ws.Range("Q5").Value = 1.648012530925445
$ws.Range("R5").Value = 14.832112778329
$ws.Range("S5").Value = 0.008238498141910647
$ws.Range("T5").Value = 0.008238498141910647

$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Lta"
$ws.Range("C6").Value = "Tnfrsf1b"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6836903333333334
$ws.Range("H6").Value = 2.051071
$ws.Range("I6").Value = 0.7427603292077943
$ws.Range("J6").Value = 0.7427603292077943
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 83.95738966666666
$ws.Range("N6").Value = 251.872169
$ws.Range("O6").Value = 0.386328906947963
$ws.Range("P6").Value = 0.386328906947963
$ws.Range("Q6").Value = 57.40085572699989
$ws.Range("R6").Value = 516.6077015429991
$ws.Range("S6").Value = 0.2869497861071564
$ws.Range("T6").Value = 0.2869497861071564

$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Lta"
$ws.Range("C7").Value = "Tnfrsf1b"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.174178
$ws.Range("H7").Value = 0.5225340000000001
$ws.Range("I7").Value = 0.1892267629264251
$ws.Range("J7").Value = 0.1892267629264251
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 18.95273633333333
$ws.Range("N7").Value = 56.858209
$ws.Range("O7").Value = 0.08721078561875104
$ws.Range("P7").Value = 0.08721078561875105
$ws.Range("Q7").Value = 3.301149709067333
$ws.Range("R7").Value = 29.710347381606
$ws.Range("S7").Value = 0.01650261465490669
$ws.Range("T7").Value = 0.01650261465490669

$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Lta"
$ws.Range("C8").Value = "Tnfrsf1b"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.174178
$ws.Range("H8").Value = 0.5225340000000001
$ws.Range("I8").Value = 0.1892267629264251
$ws.Range("J8").Value = 0.1892267629264251
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 10.54486
$ws.Range("N8").Value = 31.63458
$ws.Range("O8").Value = 0.04852204497892696
$ws.Range("P8").Value = 0.04852204497892696
$ws.Range("Q8").Value = 1.83668262508
$ws.Range("R8").Value = 16.53014362572
$ws.Range("S8").Value = 0.009181669501932747
$ws.Range("T8").Value = 0.009181669501932745

$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Lta"
$ws.Range("C9").Value = "Tnfrsf1b"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.174178
$ws.Range("H9").Value = 0.5225340000000001
$ws.Range("I9").Value = 0.1892267629264251
$ws.Range("J9").Value = 0.1892267629264251
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 101.4555613333333
$ws.Range("N9").Value = 304.366684
$ws.Range("O9").Value = 0.4668465309523581
$ws.Range("P9").Value = 0.4668465309523581
$ws.Range("Q9").Value = 17.67132676191734
$ws.Range("R9").Value = 159.041940857256
$ws.Range("S9").Value = 0.08833985783554585
$ws.Range("T9").Value = 0.08833985783554583

$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Lta"
$ws.Range("C10").Value = "Tnfrsf1b"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.174178
$ws.Range("H10").Value = 0.5225340000000001
$ws.Range("I10").Value = 0.1892267629264251
$ws.Range("J10").Value = 0.1892267629264251
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.410466333333333
$ws.Range("N10").Value = 7.231399000000001
$ws.Range("O10").Value = 0.01109173150200089
$ws.Range("P10").Value = 0.01109173150200089
$ws.Range("Q10").Value = 0.4198502050073334
$ws.Range("R10").Value = 3.778651845066001
$ws.Range("S10").Value = 0.002098852447372684
$ws.Range("T10").Value = 0.002098852447372684

$ws.Range("A11").Value = "Inflammatory-Mac"
$ws.Range("B11").Value = "Lta"
$ws.Range("C11").Value = "Tnfrsf1b"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.174178
$ws.Range("H11").Value = 0.5225340000000001
$ws.Range("I11").Value = 0.1892267629264251
$ws.Range("J11").Value = 0.1892267629264251
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 83.95738966666666
$ws.Range("N11").Value = 251.872169
$ws.Range("O11").Value = 0.386328906947963
$ws.Range("P11").Value = 0.386328906947963
$ws.Range("Q11").Value = 14.62353021736067
$ws.Range("R11").Value = 131.611771956246
$ws.Range("S11").Value = 0.07310376848666714
$ws.Range("T11").Value = 0.07310376848666712

$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Lta"
$ws.Range("C12").Value = "Tnfrsf1b"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.06260399999999999
$ws.Range("H12").Value = 0.187812
$ws.Range("I12").Value = 0.06801290786578049
$ws.Range("J12").Value = 0.06801290786578047
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 18.95273633333333
$ws.Range("N12").Value = 56.858209
$ws.Range("O12").Value = 0.08721078561875104
$ws.Range("P12").Value = 0.08721078561875105
$ws.Range("Q12").Value = 1.186517105412
$ws.Range("R12").Value = 10.678653948708
$ws.Range("S12").Value = 0.005931459127190448
$ws.Range("T12").Value = 0.005931459127190447

$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Lta"
$ws.Range("C13").Value = "Tnfrsf1b"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.06260399999999999
$ws.Range("H13").Value = 0.187812
$ws.Range("I13").Value = 0.06801290786578049
$ws.Range("J13").Value = 0.06801290786578047
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 10.54486
$ws.Range("N13").Value = 31.63458
$ws.Range("O13").Value = 0.04852204497892696
$ws.Range("P13").Value = 0.04852204497892696
$ws.Range("Q13").Value = 0.66015041544
$ws.Range("R13").Value = 5.941353738959999
$ws.Range("S13").Value = 0.003300125374611016
$ws.Range("T13").Value = 0.003300125374611015

$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Lta"
$ws.Range("C14").Value = "Tnfrsf1b"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.06260399999999999
$ws.Range("H14").Value = 0.187812
$ws.Range("I14").Value = 0.06801290786578049
$ws.Range("J14").Value = 0.06801290786578047
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 101.4555613333333
$ws.Range("N14").Value = 304.366684
$ws.Range("O14").Value = 0.4668465309523581
$ws.Range("P14").Value = 0.4668465309523581
$ws.Range("Q14").Value = 6.351523961711999
$ws.Range("R14").Value = 57.163715655408
$ws.Range("S14").Value = 0.03175159009712197
$ws.Range("T14").Value = 0.03175159009712197

$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Lta"
$ws.Range("C15").Value = "Tnfrsf1b"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.06260399999999999
$ws.Range("H15").Value = 0.187812
$ws.Range("I15").Value = 0.06801290786578049
$ws.Range("J15").Value = 0.06801290786578047
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.410466333333333
$ws.Range("N15").Value = 7.231399000000001
$ws.Range("O15").Value = 0.01109173150200089
$ws.Range("P15").Value = 0.01109173150200089
$ws.Range("Q15").Value = 0.150904834332
$ws.Range("R15").Value = 1.358143508988
$ws.Range("S15").Value = 0.0007543809127175618
$ws.Range("T15").Value = 0.0007543809127175617

$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Lta"
$ws.Range("C16").Value = "Tnfrsf1b"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.06260399999999999
$ws.Range("H16").Value = 0.187812
$ws.Range("I16").Value = 0.06801290786578049
$ws.Range("J16").Value = 0.06801290786578047
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 83.95738966666666
$ws.Range("N16").Value = 251.872169
$ws.Range("O16").Value = 0.386328906947963
$ws.Range("P16").Value = 0.386328906947963
$ws.Range("Q16").Value = 5.256068422691999
$ws.Range("R16").Value = 47.30461580422799
$ws.Range("S16").Value = 0.02627535235413949
$ws.Range("T16").Value = 0.02627535235413949

